$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.99000000000062
$ws.Range("G2").Value = [double]"4.662936703425657e-15"
$ws.Range("H2").Value = [double]"1.247229270217197e-13"
$ws.Range("K2").Value = 50.36365288928197
$ws.Range("L2").Value = "[36.299871942737084, 64.42743383582686]"
$ws.Range("M2").Value = [double]"3.013478355740062e-11"
$ws.Range("N2").Value = [double]"6.026956711480125e-11"
$ws.Range("O2").Value = 1.855395060678657
$ws.Range("P2").Value = "[1.5660792207084269, 2.144710900648888]"
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 68.79911129204837
$ws.Range("T2").Value = "[60.857344997848905, 76.74087758624782]"
$ws.Range("W2").Value = 18.31527527527571
$ws.Range("X2").Value = 17.11853853853895
$ws.Range("Y2").Value = 19.51201201201248

# Row 3
$ws.Range("E3").Value = 23.09000000000017
$ws.Range("G3").Value = [double]"5.15428255631889e-10"
$ws.Range("H3").Value = [double]"1.507550540858838e-09"
$ws.Range("K3").Value = 48.29417095871784
$ws.Range("L3").Value = "[30.195700926342425, 66.39264099109326]"
$ws.Range("M3").Value = [double]"4.884700293139588e-07"
$ws.Range("N3").Value = [double]"4.884700293139588e-07"
$ws.Range("O3").Value = -0.5534737808126167
$ws.Range("P3").Value = "[-0.9308422677303092, -0.17610529389492413]"
$ws.Range("Q3").Value = 0.004336540958253199
$ws.Range("R3").Value = 0.004336540958253199
$ws.Range("S3").Value = 66.76014308216446
$ws.Range("T3").Value = "[57.19049260794429, 76.32979355638463]"
$ws.Range("W3").Value = 2.033953953953972
$ws.Range("X3").Value = 0.6471671671671764
$ws.Range("Y3").Value = 3.420740740740768
